# "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
#
# The worker table (rows 16-17) held one employee (YESICA PAOLA GUETO BARRIOS)
# with two overdue periods (2110, 2111). The edit:
#   - adds a brand-new employee (MARLON ENRIQUE LEON ANTEQUERA, period 2207)
#     as the new first data row,
#   - keeps YESICA's two periods (now one row lower each),
#   - updates F/G amounts for the first data row,
#   - refreshes the VALOR MORA total and the trabajadores/periodos counters.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for a 3rd data row: insert a blank row at 17, shifting the old
# row 17 (and everything below, incl. the signature block) down by one.
$ws.Rows("17:17").Insert()

# The freshly inserted row 17 has no formatting yet; give it the same
# "interior" table-row style as row 16 (borders/number formats) without
# touching row 16's values.
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J17").PasteSpecial(-4122)

# Start the 3 data rows from a clean slate and retype them in the new order:
# new employee first, then the existing employee's two periods.
$ws.Range("B16:J18").ClearContents()

# Row 16: new employee MARLON ENRIQUE LEON ANTEQUERA, period 2207
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1149446282"
$ws.Range("D16").Value = "MARLON ENRIQUE LEON ANTEQUERA"
$ws.Range("E16").Value = "2207"
$ws.Range("F16").Value = 40000
$ws.Range("G16").Value = 1000000

# Row 17: existing employee, period 2111
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1007412417"
$ws.Range("D17").Value = "YESICA PAOLA GUETO BARRIOS"
$ws.Range("E17").Value = "2111"
$ws.Range("F17").Value = 72682
$ws.Range("G17").Value = 1817052

# Row 18: existing employee, period 2110
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1007412417"
$ws.Range("D18").Value = "YESICA PAOLA GUETO BARRIOS"
$ws.Range("E18").Value = "2110"
$ws.Range("F18").Value = 36341
$ws.Range("G18").Value = 1817052

# Refresh the summary figures above the table.
$ws.Range("E11").Value = 149023
$ws.Range("C13").Value = 2
$ws.Range("F13").Value = 3
